# Review data was re-keyed: each state's rows (19-89 on Sheet1) effectively
# shifted by one position relative to last week, plus a handful of ASM
# reassignments. Apply the final SS/ASM/Password values cell-by-cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(19, 2).Value = 'Punjab'
$ws.Cells.Item(19, 3).Value = 'Anand Sales|ACUSDA0820'
$ws.Cells.Item(19, 4).Value = 'Neeraj Kapoor'
$ws.Cells.Item(20, 3).Value = 'Vishal Traders|ACUSDV0085'
$ws.Cells.Item(20, 4).Value = 'Jasbeer Verma'
$ws.Cells.Item(21, 3).Value = 'Ashish & Co.|ACUSDA1159'
$ws.Cells.Item(21, 4).Value = 'Badal Srivastava'
$ws.Cells.Item(22, 3).Value = 'P.R.Sales Agency|ACUSDP0142'
$ws.Cells.Item(22, 4).Value = 'Rakesh Kumar'
$ws.Cells.Item(23, 3).Value = 'Swastik Enterprises|ACUSDS1206'
$ws.Cells.Item(23, 4).Value = 'Jasbeer Verma'
$ws.Cells.Item(24, 3).Value = 'Dalal Trading Co.|ACUSDD1126'
$ws.Cells.Item(24, 4).Value = 'Yogesh'
$ws.Cells.Item(25, 3).Value = 'Shri Radhey Krishna Trading Co.|ACUSDS1204'
$ws.Cells.Item(25, 4).Value = 'Jasbeer Verma'
$ws.Cells.Item(26, 3).Value = 'Vishal Trading Co.|ACUSDV0079'
$ws.Cells.Item(26, 4).Value = 'Pradeep Kumar'
$ws.Cells.Item(27, 3).Value = 'Shubham Overseas|ACUSDS1407'
$ws.Cells.Item(27, 4).Value = 'Yogesh'
$ws.Cells.Item(28, 3).Value = 'Rameshwarm Enterprises|ACUSDR4921'
$ws.Cells.Item(28, 4).Value = 'Ajit Pal'
$ws.Cells.Item(29, 3).Value = 'Pulkit Trading Co.|ACUSDP0132'
$ws.Cells.Item(29, 4).Value = 'Badal Srivastava'
$ws.Cells.Item(30, 3).Value = 'Paras Kumar Satish Kr|ACUSDP0929'
$ws.Cells.Item(30, 4).Value = 'Rakesh Kumar'
$ws.Cells.Item(31, 3).Value = 'Murli Wala Sales|ACUSDM0200'
$ws.Cells.Item(31, 4).Value = 'Mahesh'
$ws.Cells.Item(32, 3).Value = 'Kalyani Trading Co.|ACUSDK0087'
$ws.Cells.Item(32, 4).Value = 'Ajit Pal'
$ws.Cells.Item(33, 3).Value = 'Harshit Trading Com|ACUSDY0005'
$ws.Cells.Item(33, 4).Value = 'Badal Srivastava'
$ws.Cells.Item(34, 3).Value = 'Ravindra Associates|ACUSDR0520'
$ws.Cells.Item(34, 4).Value = 'Yogesh'
$ws.Cells.Item(35, 3).Value = 'Shakshi Enterprises|ACUSDS0338'
$ws.Cells.Item(35, 4).Value = 'Jasbeer Verma'
$ws.Cells.Item(36, 3).Value = 'Jai Shri Shyam Traders|LCLSOJ0054'
$ws.Cells.Item(37, 2).Value = 'Haryana'
$ws.Cells.Item(37, 3).Value = 'Maruti Enterprises|ACUSDM0088'
$ws.Cells.Item(37, 4).Value = 'Yogesh'
$ws.Cells.Item(38, 3).Value = 'Deepak Enterprises|ACUSDD0616'
$ws.Cells.Item(39, 3).Value = 'Katna Sales|ACUSDK0086'
$ws.Cells.Item(40, 3).Value = 'Rattanchand And Sons|ACUSDR0130'
$ws.Cells.Item(41, 3).Value = 'Amar Nath Harishchand|ACUSDA0133'
$ws.Cells.Item(42, 2).Value = 'Himachal Pradesh'
$ws.Cells.Item(42, 3).Value = 'Shashank Enteprises|ACUSDS4809'
$ws.Cells.Item(42, 4).Value = 'Parmod Kumar'
$ws.Cells.Item(43, 3).Value = 'Vinay Traders|ACUSDV0207'
$ws.Cells.Item(43, 4).Value = 'Vinay Pandey'
$ws.Cells.Item(49, 3).Value = 'Anant Wardrobe|ACUSD00073'
$ws.Cells.Item(49, 4).Value = 'Brajesh Sharma'
$ws.Cells.Item(50, 3).Value = 'Aakash Trading Co|ACUSDA0830'
$ws.Cells.Item(50, 4).Value = 'Ankit Bhardwaj'
$ws.Cells.Item(51, 3).Value = 'Churaman Biharilal Asati|ACUSDC0531'
$ws.Cells.Item(51, 4).Value = 'Anand Gupta'
$ws.Cells.Item(52, 3).Value = 'Varun Enterprises|ACUSDV0208'
$ws.Cells.Item(52, 4).Value = 'Vinay Pandey'
$ws.Cells.Item(53, 3).Value = 'Amit Sales|ACUSDA0137'
$ws.Cells.Item(53, 4).Value = 'Anand Gupta'
$ws.Cells.Item(54, 3).Value = 'Shanti Marketing|ACUSDS4808'
$ws.Cells.Item(54, 4).Value = 'Vinay Pandey'
$ws.Cells.Item(55, 3).Value = 'Abhinandan Enterprises|ACUSDA0088'
$ws.Cells.Item(55, 4).Value = 'Brajesh Sharma'
$ws.Cells.Item(56, 3).Value = 'O.S. Agency|ACUSDO0028'
$ws.Cells.Item(57, 3).Value = 'R.D. Plastic|ACUSDR0127'
$ws.Cells.Item(57, 4).Value = 'Vinay Pandey'
$ws.Cells.Item(58, 3).Value = 'Krishna Enterprises|'
$ws.Cells.Item(59, 3).Value = 'Aggarwal Agency|ACUSDA0410'
$ws.Cells.Item(59, 4).Value = 'Ankit Bhardwaj'
$ws.Cells.Item(60, 4).Value = 'Vinay Pandey'
$ws.Cells.Item(61, 3).Value = 'Bhagwan Das Kirana And General Store|'
$ws.Cells.Item(61, 4).Value = $null
$ws.Cells.Item(62, 2).Value = 'Madhya Pradesh'
$ws.Cells.Item(62, 3).Value = 'Reliable Industry|ACUSDR0204'
$ws.Cells.Item(62, 4).Value = 'Vinay Pandey'
$ws.Cells.Item(63, 3).Value = 'Prabhat Agency|ACUSDP2385'
$ws.Cells.Item(63, 4).Value = 'Rakesh Tripathi'
$ws.Cells.Item(63, 5).Value = 2025
$ws.Cells.Item(64, 3).Value = 'Prabhat Trading Com.|ACUSDP0183'
$ws.Cells.Item(64, 4).Value = 'Anil Jamadar Singh'
$ws.Cells.Item(64, 5).Value = 5808
$ws.Cells.Item(65, 3).Value = 'Adhira Agency|ACUSDSS005'
$ws.Cells.Item(65, 4).Value = 'Vacant'
$ws.Cells.Item(65, 5).Value = 1164
$ws.Cells.Item(66, 3).Value = 'Optimize Crusder|ACUSDO0026'
$ws.Cells.Item(66, 4).Value = 'Kapil Sharma'
$ws.Cells.Item(66, 5).Value = 4920
$ws.Cells.Item(67, 3).Value = 'Jainil Traders|ACUSDJ0023'
$ws.Cells.Item(68, 3).Value = 'Krishna Agencies|ACUSDSS003'
$ws.Cells.Item(68, 4).Value = 'Vacant'
$ws.Cells.Item(68, 5).Value = 1164
$ws.Cells.Item(69, 3).Value = 'Bharat Traders|ACUSDB0039'
$ws.Cells.Item(69, 4).Value = 'Kapil Sharma'
$ws.Cells.Item(69, 5).Value = 4920
$ws.Cells.Item(70, 3).Value = 'J.K. Brothers|ACUSDJ4838'
$ws.Cells.Item(71, 3).Value = 'Jai Mata Di Agency|ACUSDSS001'
$ws.Cells.Item(71, 4).Value = 'Vacant'
$ws.Cells.Item(71, 5).Value = 1164
$ws.Cells.Item(72, 3).Value = 'Sunrise Enterprises|ACUSDS0771'
$ws.Cells.Item(72, 4).Value = 'Dinesh Sharma'
$ws.Cells.Item(72, 5).Value = 7757
$ws.Cells.Item(73, 3).Value = 'Nidhi Sales|ACUSDN4807'
$ws.Cells.Item(73, 4).Value = 'Anil Jamadar Singh'
$ws.Cells.Item(73, 5).Value = 5808
$ws.Cells.Item(74, 3).Value = 'Krishna Marketing|ACUSDSS006'
$ws.Cells.Item(74, 4).Value = 'Vacant'
$ws.Cells.Item(74, 5).Value = 1164
$ws.Cells.Item(75, 3).Value = 'Ajinkya Traders|ACUSDA0146'
$ws.Cells.Item(75, 4).Value = 'Shailesh Surve'
$ws.Cells.Item(75, 5).Value = 7136
$ws.Cells.Item(76, 3).Value = 'Tirupati Agencies|ACUSDSS002'
$ws.Cells.Item(77, 3).Value = 'Shree Datta Agency|ACUSDSS004'
$ws.Cells.Item(77, 4).Value = 'Rakesh Tripathi'
$ws.Cells.Item(77, 5).Value = 2025
$ws.Cells.Item(78, 2).Value = 'Maharashtra'
$ws.Cells.Item(78, 3).Value = 'Navinya Enterprises|ACUSDN0052'
$ws.Cells.Item(78, 4).Value = 'Shailesh Surve'
$ws.Cells.Item(78, 5).Value = 7136
$ws.Cells.Item(79, 3).Value = 'Garg Sales Corporation|ACUSDG1091'
$ws.Cells.Item(79, 4).Value = 'Neeraj Kapoor'
$ws.Cells.Item(80, 3).Value = 'Bajarang Sweet House (Conf)|ACUSDB0963'
$ws.Cells.Item(81, 3).Value = 'Sant Enterprises|ACUSDS1100'
$ws.Cells.Item(81, 4).Value = 'Jasbir Verma'
$ws.Cells.Item(82, 3).Value = 'Shivam Agency|ACUSDSS896'
$ws.Cells.Item(82, 4).Value = 'Pradeep'
$ws.Cells.Item(83, 3).Value = 'Jeet Ram And Sons|ACUSD00063'
$ws.Cells.Item(83, 4).Value = 'Neeraj Kapoor'
$ws.Cells.Item(84, 3).Value = 'Commander Trade Link|ACUSDC1514'
$ws.Cells.Item(84, 4).Value = 'Jasbir Verma'
$ws.Cells.Item(85, 3).Value = 'Saurabh Traders|ACUSDSS106'
$ws.Cells.Item(85, 4).Value = 'Pradeep'
$ws.Cells.Item(86, 3).Value = 'Suri Enterprises|ACUSD00062'
$ws.Cells.Item(86, 4).Value = 'Pradeep'
$ws.Cells.Item(87, 3).Value = 'Grace Drinks Pvt Ltd|ACUSDG0174'
$ws.Cells.Item(87, 4).Value = 'Rohit Medirata'
$ws.Cells.Item(88, 3).Value = 'Bhana Ram & Sons|ACUSDB0473'
$ws.Cells.Item(89, 3).Value = 'Grace Drinks Pvt Ltd|ACUSDG0174'
